# Update the Daily Orders sheet: order status changes from READY to DELIVERED
$wb = $excel.ActiveWorkbook

$ordersWs = $wb.Worksheets.Item("Daily Orders")
$ordersWs.Range("H2").Value = "DELIVERED"

# Reflect the status change in the Summary sheet counts:
# the "Ready" count decreases by one, the "Delivered" count increases by one.
$summaryWs = $wb.Worksheets.Item("Summary")
$summaryWs.Range("D2").Value = 0
$summaryWs.Range("E2").Value = 1
